# Added Hit Rate Evaluation
#
# On the "Web Interface" sheet's sub-task table (Table6), insert a new
# task row before the existing "Add a second page showing the evaluation
# results" row, describing the fix for the Collaborative Recommenders,
# and lower the priority of the "Allow selectable evaluation options"
# task (which shifts down to the end of the table) to Medium.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Web Interface")
$lo = $ws.ListObjects.Item(1)

# Shift the existing rows 8-9 down to 9-10, keeping their values/styles,
# and make room for the new row at position 8.
$ws.Rows.Item(8).Insert()

# Grow the table to include the newly inserted row.
$lo.Resize($ws.Range("A1:D10"))

# Fill in the new task row.
$ws.Range("A8").Value = "Fix Collaborative Recommenders so they show their recommendations results"
$ws.Range("B8").Value = "High"
$ws.Range("C8").Value = "No"

# Match the "Bad" (red) styling used by the other unfinished rows, sized
# like the rest of the manually-typed rows in this table.
$ws.Range("B8").Style = "Bad"
$ws.Range("B8").Font.Size = 11
$ws.Range("C8").Style = "Bad"
$ws.Range("C8").Font.Size = 11

# This task's text wraps onto three lines, same as the other long entries.
$ws.Rows.Item(8).RowHeight = 47.25

# The "Allow selectable evaluation options" task (now row 10) is
# downgraded from High to Medium priority.
$ws.Range("B10").Value = "Medium"
$ws.Range("B10").Style = "Neutral"
$ws.Range("B10").Font.Size = 11

# Move the active selection, matching where the author ended up editing.
$ws.Range("I5").Select()
